# Applies the cryptos-list refresh described in the commit:
#  - OKB is inserted as a new row; every coin from "Dogecoin" down to
#    "NEARProtocol" shifts down one row (Cronos falls off the bottom).
#  - Price (column D) and Volume(1h) (column E) are refreshed for every row.
$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values that Excel would otherwise
# auto-convert to a number (e.g. "1.001", "45.39") are prefixed with a
# leading apostrophe so they are stored as literal text, matching the
# workbook (the Price column is text, not numeric).
$updates = [ordered]@{
    'D2' = '25.748.12'
    'E2' = '  -4.24%  '
    'D3' = '1.812.37'
    'E3' = '  -3.23%  '
    'D5' = '''276.37'
    'E5' = '  -8.42%  '
    'D6' = '''1.001'
    'E6' = '  -0.06%  '
    'D7' = '''0.5058'
    'E7' = '  -4.57%  '
    'E8' = '  -6.51%  '
    'B9' = 'OKB'
    'C9' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D9' = '''45.39'
    'E9' = '  -0.20%  '
    'B10' = 'Dogecoin'
    'C10' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'D10' = '''0.06670'
    'E10' = '  -6.90%  '
    'B11' = 'Solana'
    'C11' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    'D11' = '''20.00'
    'E11' = '  -7.42%  '
    'B12' = 'Polygon'
    'C12' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D12' = '''0.8305'
    'E12' = '  -6.14%  '
    'B13' = 'TRON'
    'C13' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D13' = '''0.07894'
    'E13' = '  -2.93%  '
    'B14' = 'WrappedEther'
    'C14' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D14' = '1.819.05'
    'E14' = '  -2.59%  '
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D15' = '''5.072'
    'E15' = '  -3.81%  '
    'B16' = 'Litecoin'
    'C16' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D16' = '''87.44'
    'E16' = '  -6.00%  '
    'B17' = 'BinanceUSD'
    'C17' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D17' = '''1.000'
    'E17' = '  -0.15%  '
    'B18' = 'Avalanche'
    'C18' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D18' = '''14.00'
    'E18' = '  -4.81%  '
    'B19' = 'ShibaInu'
    'C19' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D19' = '''0.000008027'
    'E19' = '  -5.86%  '
    'B20' = 'Dai'
    'C20' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D20' = '''1.000'
    'E20' = '  +0.02%  '
    'B21' = 'WrappedBTC'
    'C21' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D21' = '25.790.41'
    'E21' = '  -4.13%  '
    'B22' = 'Uniswap'
    'C22' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D22' = '''4.715'
    'E22' = '  -5.22%  '
    'B23' = 'Cosmos'
    'C23' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D23' = '''9.998'
    'E23' = '  -6.44%  '
    'B24' = 'Chainlink'
    'C24' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D24' = '''6.044'
    'E24' = '  -5.05%  '
    'B25' = 'Monero'
    'C25' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D25' = '''142.18'
    'E25' = '  -3.54%  '
    'B26' = 'LidoDAOToken'
    'C26' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D26' = '''2.177'
    'E26' = '  -3.89%  '
    'B27' = 'Toncoin'
    'C27' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D27' = '''1.667'
    'E27' = '  -3.87%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D28' = '''17.02'
    'E28' = '  -5.54%  '
    'B29' = 'BitcoinCash'
    'C29' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D29' = '''109.41'
    'E29' = '  -4.45%  '
    'B30' = 'InternetComputer(DFINITY)'
    'C30' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D30' = '''4.318'
    'E30' = '  -8.85%  '
    'B31' = 'Filecoin'
    'C31' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D31' = '''4.222'
    'E31' = '  -7.68%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D32' = '''0.08792'
    'E32' = '  -3.41%  '
    'B33' = 'Hedera'
    'C33' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D33' = '''0.04859'
    'E33' = '  -2.39%  '
    'B34' = 'ImmutableX'
    'C34' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D34' = '''0.7256'
    'E34' = '  -8.84%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D35' = '''1.135'
    'E35' = '  -3.18%  '
    'B36' = 'HuobiToken'
    'C36' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D36' = '''2.874'
    'E36' = '  -3.80%  '
    'D37' = '''3.152'
    'E37' = '  -1.34%  '
    'B38' = 'Frax'
    'C38' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'D38' = '''1.0000'
    'E38' = '  +0.01%  '
    'B39' = 'TheSandbox'
    'C39' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D39' = '''0.5190'
    'E39' = '  -11.22%  '
    'B40' = 'VeChain'
    'C40' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D40' = '''0.01839'
    'E40' = '  -5.63%  '
    'B41' = 'RenderToken'
    'C41' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D41' = '''2.263'
    'E41' = '  -13.45%  '
    'B42' = 'TrustWalletToken'
    'C42' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D42' = '''0.9509'
    'E42' = '  -11.16%  '
    'B43' = 'Quant'
    'C43' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D43' = '''112.95'
    'E43' = '  -2.28%  '
    'B44' = 'FraxShare'
    'C44' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D44' = '''6.129'
    'E44' = '  -7.19%  '
    'B45' = 'Aptos'
    'C45' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D45' = '''8.055'
    'E45' = '  -9.25%  '
    'B46' = 'PaxDollar'
    'C46' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D46' = '''1.0000'
    'E46' = '  -0.02%  '
    'B47' = 'Decentraland'
    'C47' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D47' = '''0.4553'
    'E47' = '  -9.83%  '
    'B48' = 'Algorand'
    'C48' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D48' = '''0.1356'
    'E48' = '  -9.18%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D49' = '''9.286'
    'E49' = '  -7.02%  '
    'B50' = 'Elrond'
    'C50' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D50' = '''36.27'
    'E50' = '  -4.48%  '
    'B51' = 'NEARProtocol'
    'C51' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D51' = '''1.499'
    'E51' = '  -6.91%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
